$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (Late/Outstanding columns shift
# right: old N->O, O->P, P->Q). This mirrors a manual "Insert Column" on the
# repayment schedule sheet to make room for an extra (variable instalment)
# column.
$ws.Columns("N").Insert() | Out-Null

# The freshly inserted column inherits the width of the column immediately
# to its left (column M), matching Excel's own Insert-Column behaviour.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make the "Repayment schedule" sheet the active tab (it becomes the
# selected/visible sheet after the edit).
$ws.Activate() | Out-Null

# Move the selection to just past the new data (one column/row beyond the
# inserted column), matching where the cursor ends up after the edit.
$ws.Range("R5").Select() | Out-Null
